$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh inserts two new observation rows (new price points
# dated 44449) right before the existing row that used to be 191, pushing the
# rest of the "Cilantro" block (old rows 191-208) down by two (new rows
# 193-210).
$ws.Rows.Item(191).Insert()
$ws.Rows.Item(191).Insert()

# --- New row 191 ---
$ws.Cells.Item(191,1).Value2  = 10
$ws.Cells.Item(191,2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(191,3).Value2  = "La Araucanía"
$ws.Cells.Item(191,4).Value2  = 44449
$ws.Cells.Item(191,5).Value2  = 9
$ws.Cells.Item(191,6).Value2  = 100112040
$ws.Cells.Item(191,7).Value2  = "Cilantro"
$ws.Cells.Item(191,8).Value2  = "Sin especificar"
$ws.Cells.Item(191,9).Value2  = "Primera"
$ws.Cells.Item(191,10).Value2 = 125
$ws.Cells.Item(191,11).Value2 = 4000
$ws.Cells.Item(191,12).Value2 = 4000
$ws.Cells.Item(191,13).Value2 = 4000
$ws.Cells.Item(191,14).Value2 = "$/docena de atados (2 kilos)"
$ws.Cells.Item(191,15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(191,16).Value2 = 2000
$ws.Cells.Item(191,17).Value2 = 2
$ws.Cells.Item(191,18).Value2 = "Hortaliza"

# --- New row 192 ---
$ws.Cells.Item(192,1).Value2  = 10
$ws.Cells.Item(192,2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(192,3).Value2  = "La Araucanía"
$ws.Cells.Item(192,4).Value2  = 44449
$ws.Cells.Item(192,5).Value2  = 9
$ws.Cells.Item(192,6).Value2  = 100112040
$ws.Cells.Item(192,7).Value2  = "Cilantro"
$ws.Cells.Item(192,8).Value2  = "Sin especificar"
$ws.Cells.Item(192,9).Value2  = "Primera"
$ws.Cells.Item(192,10).Value2 = 75
$ws.Cells.Item(192,11).Value2 = 4500
$ws.Cells.Item(192,12).Value2 = 4500
$ws.Cells.Item(192,13).Value2 = 4500
$ws.Cells.Item(192,14).Value2 = "$/docena de atados (2 kilos)"
$ws.Cells.Item(192,15).Value2 = "Región Metropolitana"
$ws.Cells.Item(192,16).Value2 = 2250
$ws.Cells.Item(192,17).Value2 = 2
$ws.Cells.Item(192,18).Value2 = "Hortaliza"
